$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 315, shifting rows 315:357 down to 316:358.
$ws.Rows("315:315").Insert(1)

# Populate the new row 315 with the new "weekly" record.
$ws.Cells.Item(315, 1).Value = 7
$ws.Cells.Item(315, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(315, 3).Value = "Ñuble"
$ws.Cells.Item(315, 4).Value = 45131
$ws.Cells.Item(315, 5).Value = 16
$ws.Cells.Item(315, 6).Value = 100112043
$ws.Cells.Item(315, 7).Value = "Pepino ensalada"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 120
$ws.Cells.Item(315, 11).Value = 15000
$ws.Cells.Item(315, 12).Value = 17000
$ws.Cells.Item(315, 13).Value = 15667
$ws.Cells.Item(315, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(315, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(315, 16).Value = 261
$ws.Cells.Item(315, 17).Value = 60
$ws.Cells.Item(315, 18).Value = "Hortaliza"
